$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "month" formulas in Q2 and Q3 so they look at the date 30 days ago
# instead of today's date (fix for Agenda Mobile/Portal view).
$ws.Range("Q2").Formula = '=TEXT(TODAY()-30,"mmmm")'
$ws.Range("Q3").Formula = '=TEXT(TODAY()-30,"mmmm")'

# Update the active selection to Q4, matching the saved workbook state.
$ws.Range("Q4").Select()
